$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 322
$ws1.Range("F4").Value = 1299
$ws1.Range("F5").Value = 641

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 322
$ws4.Range("F4").Value = 1299
$ws4.Range("F6").Value = 641
